$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-04-14 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-15 Monday", 2)

# Update the multiplication answers in the table, cell by cell, so that
# duplicate source values (e.g. "498x6=2988" appearing twice) are each
# replaced with their own distinct target value instead of both being
# overwritten identically.
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "852×3=2556"
$tbl.Cell(1, 2).Range.Text = "298×2=596"
$tbl.Cell(1, 3).Range.Text = "154×5=770"
$tbl.Cell(1, 4).Range.Text = "796×3=2388"
$tbl.Cell(1, 5).Range.Text = "606×6=3636"

$tbl.Cell(5, 1).Range.Text = "475×8=3800"
$tbl.Cell(5, 2).Range.Text = "945×8=7560"
$tbl.Cell(5, 3).Range.Text = "844×8=6752"
$tbl.Cell(5, 4).Range.Text = "289×7=2023"
$tbl.Cell(5, 5).Range.Text = "332×9=2988"

$tbl.Cell(10, 1).Range.Text = "320×5=1600"
$tbl.Cell(10, 2).Range.Text = "830×4=3320"
$tbl.Cell(10, 3).Range.Text = "354×9=3186"
$tbl.Cell(10, 4).Range.Text = "563×8=4504"
$tbl.Cell(10, 5).Range.Text = "788×6=4728"

$tbl.Cell(15, 1).Range.Text = "130×4=520"
$tbl.Cell(15, 2).Range.Text = "271×2=542"
$tbl.Cell(15, 3).Range.Text = "728×8=5824"
$tbl.Cell(15, 4).Range.Text = "633×3=1899"
$tbl.Cell(15, 5).Range.Text = "891×9=8019"

$tbl.Cell(20, 1).Range.Text = "430×8=3440"
$tbl.Cell(20, 2).Range.Text = "647×6=3882"
$tbl.Cell(20, 3).Range.Text = "756×3=2268"
$tbl.Cell(20, 4).Range.Text = "400×7=2800"
$tbl.Cell(20, 5).Range.Text = "648×4=2592"
